# Daily data-refresh upload: bump the USD Amount figure in cell T2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

$ws.Range("T2").Value = 552005
